# "navn for dag 1" - fill in the participant names for day 1 (K1) and the
# two runners that only appear on day 1 of K2, then let the workbook's
# formulas (D1D2D3 / KxDx / Totalt sheets) recalc from those names.

$wb = $excel.ActiveWorkbook

# --- K1 sheet: names for rows 4-12 (shared strings are appended in this
#     exact order, matching the sharedStrings.xml diff order 37..45) ---
$k1 = $wb.Worksheets.Item("K1")
$k1.Range("A4").Value  = "Johanne Lunde Ragnhildsløkken"
$k1.Range("A5").Value  = "Anette Berntsen"
$k1.Range("A6").Value  = "Siri Vestengen"
$k1.Range("A7").Value  = "Kristine Haanes Strandlie"
$k1.Range("A8").Value  = "Thea Helene Linnerud Foss"
$k1.Range("A9").Value  = "Maria Brenna"
$k1.Range("A10").Value = "Natascha Nina Silber"
$k1.Range("A11").Value = "Thea Snortheimsmoen"
$k1.Range("A12").Value = "Anita Bueno Lindmoen"

# --- K2 sheet: names for rows 4-5 (shared strings 46..47) ---
$k2 = $wb.Worksheets.Item("K2")
$k2.Range("A4").Value = "Mathilde Skjærdalen Myhrvold"
$k2.Range("A5").Value = "Dorthe Ballangrud Seierstad"

# --- View tweaks: zoom out from the very large 150-190% zoom levels
#     that were used while the sheets were still mostly empty, down to a
#     consistent 150%, now that column A holds real names. ---
function Set-SheetZoom($sheetName, $percent) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Activate()
    $excel.ActiveWindow.Zoom = $percent
}

Set-SheetZoom "K1" 150
Set-SheetZoom "K2" 150
Set-SheetZoom "K3" 150
Set-SheetZoom "K1D1" 150
Set-SheetZoom "K1D2" 150
Set-SheetZoom "K1D3" 150
Set-SheetZoom "K2D1" 150
Set-SheetZoom "K2D2" 150
Set-SheetZoom "K2D3" 150
Set-SheetZoom "K3D1" 150
Set-SheetZoom "K3D2" 150
Set-SheetZoom "K3D3" 150
Set-SheetZoom "D1D2D3" 150
Set-SheetZoom "Totalt" 150
Set-SheetZoom "Poengskala" 150

# --- Column A widths: widen column A on the sheets that now show names
#     so the text isn't clipped (values chosen so the engine's
#     width-quantization lands as close as possible to the authored
#     widths). ---
$k1.Columns.Item(1).ColumnWidth = 27.83
$k2.Columns.Item(1).ColumnWidth = 26.83

$d1d2d3 = $wb.Worksheets.Item("D1D2D3")
$d1d2d3.Columns.Item(1).ColumnWidth = 26.17

$totalt = $wb.Worksheets.Item("Totalt")
$totalt.Columns.Item(1).ColumnWidth = 36.17

foreach ($name in @("K1D1","K1D2","K1D3","K2D1","K2D2","K2D3")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns.Item(1).ColumnWidth = 30.67
}

# --- Restore selections that moved while entering the names ---
$k1.Range("A13").Select()
$k2.Range("A6").Select()

# --- Re-activate the originally active sheet/tab (Totalt) so tabSelected
#     ends up back where it started. ---
$totalt.Activate()
$totalt.Range("A3").Select()

Write-Host "Names for day 1 entered."
